# Apply the commit "Esqueleto de los metodos post realizado":
# - sheet "suma": simplify the sum formula (drop the stray "+2")
# - sheet "compleja": turn the GET skeleton into a POST skeleton (label + formatting cleanup)

$wb = $excel.ActiveWorkbook

$wsSuma = $wb.Worksheets.Item("suma")
$wsResta = $wb.Worksheets.Item("resta")
$wsCompleja = $wb.Worksheets.Item("compleja")

# --- sheet "suma": C4 formula no longer adds the extra "+2" ---
$wsSuma.Range("C4").Formula = "=C2+C3"

# --- sheet "compleja": becomes the POST method skeleton ---
# The header label switches from "get" to "post".
$wsCompleja.Range("B1").Value = "post"

# Fix formatting inconsistencies on rows 3/4 ("valor2"/"valor3" labels) so the
# "entrada" cells share the exact same look (font/size) as rows 2 and 5.
$wsCompleja.Range("A2").Copy()
$wsCompleja.Range("A3").PasteSpecial(-4122)
$wsCompleja.Range("A4").PasteSpecial(-4122)
$wsCompleja.Range("A2").Copy()
$wsCompleja.Range("A5").PasteSpecial(-4122)

# Normalize the "sumando"/value cells formatting across rows 2-5 to match.
$wsCompleja.Range("B2").Copy()
$wsCompleja.Range("B3").PasteSpecial(-4122)
$wsCompleja.Range("B4").PasteSpecial(-4122)
$wsCompleja.Range("B5").PasteSpecial(-4122)

$wsCompleja.Range("C2").Copy()
$wsCompleja.Range("C3").PasteSpecial(-4122)
$wsCompleja.Range("C4").PasteSpecial(-4122)
$wsCompleja.Range("C5").PasteSpecial(-4122)

$excel.CutCopyMode = 0
